$d = $word.ActiveDocument

# Locate the heading paragraph "2.6 Modelo E/R" and bump its font size
# from 14pt (sz/szCs 28) to 16pt (sz/szCs 32), on both the paragraph
# mark run properties and the text run itself.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "2.6 Modelo E/R*") {
        $r = $p.Range
        $r.Font.Size = 16
        $r.Font.SizeBi = 16
    }
}
